# The workbook has two sheets: "Planning prévisionnel" and "Planning
# effectif". The edit happened on "Planning effectif" (the active sheet),
# column H ("7e jour"), recording the hours actually spent that day on
# several tasks (filtering, user management, film deletion/blocking work).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning effectif")
$ws.Activate()

$ws.Range("H5").Value = 0.020833333333333332
$ws.Range("H9").Value = 0.020833333333333332
$ws.Range("H14").Value = 0.020833333333333332
$ws.Range("H16").Value = 0.020833333333333332
$ws.Range("H21").Value = 0.041666666666666664
$ws.Range("H27").Value = 0.020833333333333332
$ws.Range("H28").Value = 0.020833333333333332
$ws.Range("H30").Value = 0.041666666666666664
$ws.Range("H31").Value = 0.020833333333333332
$ws.Range("H38").Value = 0.083333333333333329

# The user ended up with H21 selected when they saved.
$ws.Range("H21").Select()

# The Excel window was also moved/resized on save (best effort - the
# hosted workbook-view geometry may not be reflected in every runtime).
$excel.Left = 28680
$excel.Top = -120
$excel.Width = 29040
$excel.Height = 16440
